# The deck currently uses the "Integral" design (ppt/theme/theme2.xml, the
# theme shared by the presentation + the single slide master). The edit
# swaps the active design back to the default "Office Theme" palette -
# i.e. every theme colour slot (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink)
# changes from the Integral values to the stock Office values. The font
# scheme / format scheme are identical between the two named themes, so
# only the colour scheme needs to change.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ThemeColorScheme exposes all 12 theme colour slots (unlike the legacy
# 8-slot ColorScheme) in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1..accent6, 11 hlink, 12 folHlink
$tcs = $s.ThemeColorScheme

function HexToOleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeTheme = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

foreach ($slot in $officeTheme.Keys) {
    $tcs.Colors($slot).RGB = HexToOleColor($officeTheme[$slot])
}

# Rename the design/theme to match (no-op on hosts that keep this
# read-only, harmless otherwise).
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
try { $p.Designs.Item(1).SlideMaster.Theme.Name = "Office Theme" } catch {}
